$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K (column G) values per row, regenerated from source data
$kValues = @{
    2 = 1
    3 = 1
    4 = 0
    5 = 2
    6 = 0
    7 = 0
    8 = 1
    9 = 0
    10 = 0
    11 = 1
    12 = 2
    13 = 3
    14 = 0
    15 = 1
    16 = 2
    17 = 3
    18 = 2
    19 = 1
    21 = 1
    22 = 2
    23 = 2
    24 = 2
    25 = 0
    26 = 1
    28 = 0
    29 = 0
    30 = 2
    31 = 0
    32 = 2
    33 = 2
    34 = 0
    35 = 2
    37 = 2
    38 = 0
    39 = 1
    40 = 3
    41 = 2
    42 = 2
    43 = 2
    44 = 0
    45 = 3
    46 = 1
    47 = 0
    48 = 1
    49 = 0
    50 = 2
    51 = 1
    52 = 1
    53 = 2
    54 = 1
    55 = 2
    56 = 0
    57 = 1
    59 = 3
    60 = 0
    61 = 1
    62 = 2
    63 = 2
    64 = 1
    65 = 2
    66 = 1
    67 = 1
    68 = 0
    71 = 2
    72 = 1
    73 = 2
    74 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
